$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.595.47'
$ws.Range("E2").Value = '  -0.44%  '

$ws.Range("D3").Value = '2.289.51'
$ws.Range("E3").Value = '  -1.34%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = "'96.21"
$ws.Range("E5").Value = '  +4.23%  '

$ws.Range("D6").Value = "'267.95"
$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("D7").Value = "'0.623"

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = '  -1.76%  '

$ws.Range("D10").Value = "'45.94"
$ws.Range("E10").Value = '  +2.41%  '

$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = '  -0.71%  '

$ws.Range("D12").Value = "'7.88"
$ws.Range("E12").Value = '  -1.60%  '

$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = '  +0.27%  '

$ws.Range("D14").Value = '2.632.85'
$ws.Range("E14").Value = '  -1.18%  '

$ws.Range("D15").Value = "'15.16"
$ws.Range("E15").Value = '  -0.99%  '

$ws.Range("D16").Value = "'0.851"
$ws.Range("E16").Value = '  -0.61%  '

$ws.Range("D17").Value = '2.287.73'
$ws.Range("E17").Value = '  -1.46%  '

$ws.Range("D18").Value = '43.604.52'
$ws.Range("E18").Value = '  -0.17%  '

$ws.Range("E19").Value = '  +1.35%  '

$ws.Range("E20").Value = '  -1.92%  '

$ws.Range("D21").Value = "'72.29"
$ws.Range("E21").Value = '  +1.39%  '

$ws.Range("D22").Value = "'2.50"
$ws.Range("E22").Value = '  +9.76%  '

$ws.Range("D23").Value = "'232.81"
$ws.Range("E23").Value = '  -3.66%  '

$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = '  -5.40%  '

$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").Value = '  +2.15%  '

$ws.Range("D27").Value = "'11.21"

$ws.Range("E28").Value = '  +2.81%  '

$ws.Range("D29").Value = "'40.13"
$ws.Range("E29").Value = '  +3.23%  '

$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("D31").Value = "'175.42"
$ws.Range("E31").Value = '  +1.57%  '

$ws.Range("D32").Value = "'21.85"
$ws.Range("E32").Value = '  -3.03%  '

$ws.Range("D33").Value = "'0.0893"
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").Value = "'5.37"
$ws.Range("E34").Value = '  -3.25%  '

$ws.Range("E35").Value = '  -0.50%  '

$ws.Range("E36").Value = '  -2.37%  '

$ws.Range("D38").Value = "'4.36"
$ws.Range("E38").Value = '  -3.19%  '

$ws.Range("E39").Value = '  +1.76%  '

$ws.Range("D40").Value = "'0.243"
$ws.Range("E40").Value = '  +2.59%  '

$ws.Range("E41").Value = '  -0.10%  '

$ws.Range("D42").Value = "'12.25"
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("E43").Value = '  +1.36%  '

$ws.Range("D44").Value = "'64.85"
$ws.Range("E44").Value = '  +5.89%  '

$ws.Range("E45").Value = '  -1.46%  '

$ws.Range("E46").Value = '  -4.17%  '

$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("D48").Value = "'97.36"
$ws.Range("E48").Value = '  -3.18%  '

$ws.Range("E49").Value = '  -0.51%  '

$ws.Range("D50").Value = '2.512.01'
$ws.Range("E50").Value = '  -1.22%  '

$ws.Range("E51").Value = '  +8.37%  '

